# Apply hybrid bold + color highlighting to quantitative impact metrics
# (percentages, dollar amounts, large numbers) across the resume body.
#
# Word's Font.Color takes a BGR-packed integer (0x00BBGGRR) - i.e. the
# VBA RGB(r,g,b) = r + g*256 + b*65536 encoding - so precompute the value
# for the target highlight color #2C3E50 (R=0x2C, G=0x3E, B=0x50).
$HighlightColor = 0x2C + (0x3E * 256) + (0x50 * 65536)

$d = $word.ActiveDocument

# Applies bold + the highlight color to the first occurrence of $text found
# strictly between $script:cursor and $script:paraEnd, then advances the
# cursor past the match so repeated/duplicate numbers within the same
# paragraph are highlighted left-to-right in order.
function Set-Highlight([string]$text) {
    $searchRange = $d.Range($script:cursor, $script:paraEnd)
    $found = $searchRange.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $searchRange.Font.Bold = 1
        $searchRange.Font.Color = $HighlightColor
        $script:cursor = $searchRange.End
    }
    return $found
}

# Paragraph index (1-based, Word COM convention) -> ordered list of
# substrings to bold+color within that paragraph's text.
$edits = @(
    @{ Index = 10; Segments = @("23%", "64%") },
    @{ Index = 12; Segments = @("±4.2%", "±2.1%", "71%", "87%") },
    @{ Index = 13; Segments = @("73.5%", '$4.7M') },
    @{ Index = 14; Segments = @('$2') },
    @{ Index = 34; Segments = @("57%") },
    @{ Index = 50; Segments = @("±4.2%", "±2.1%") },
    @{ Index = 51; Segments = @("71%", "87%") },
    @{ Index = 52; Segments = @("34%", "28%") }
)

foreach ($edit in $edits) {
    $p = $d.Paragraphs.Item($edit.Index)
    $script:cursor = $p.Range.Start
    $script:paraEnd = $p.Range.End
    foreach ($seg in $edit.Segments) {
        Set-Highlight $seg | Out-Null
    }
}
